$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4466.143
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 5631.5
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 5631.5
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6879.5
$ws.Range("H65").Value = 4466.143
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 5631.5
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 28157.5
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -34397.5
$ws.Range("H92").Value = 2274.3125
$ws.Range("I92").Value = 2299
$ws.Range("J92").Value = 2220
$ws.Range("K92").Value = 2299
$ws.Range("L92").Value = 2220
$ws.Range("M92").Value = -1051
$ws.Range("N92").Value = -4716
$ws.Range("H100").Value = 1888
$ws.Range("I100").Value = 1664
$ws.Range("K100").Value = 1664
$ws.Range("M100").Value = -1123
$ws.Range("H106").Value = 10494.5
$ws.Range("I106").Value = 10839.77
$ws.Range("K106").Value = 10839.77
$ws.Range("M106").Value = -10208.77
$ws.Range("H127").Value = 1085.8636
$ws.Range("I127").Value = 565.93335
$ws.Range("K127").Value = 1697.80005
$ws.Range("M127").Value = 3262.19995
$ws.Range("H132").Value = 4570549
$ws.Range("I132").Value = 5749783.5
$ws.Range("K132").Value = 17249350.5
$ws.Range("M132").Value = -17246820.5
$ws.Range("H138").Value = 1579.6632
$ws.Range("I138").Value = 885.8333
$ws.Range("K138").Value = 2657.4999
$ws.Range("M138").Value = 2482.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2883.7747
$ws.Range("I32").Value = 2641.2754
$ws.Range("J32").Value = 11250
$ws.Range("K32").Value = 2641.2754
$ws.Range("L32").Value = 11250
$ws.Range("M32").Value = -2354.2754
$ws.Range("N32").Value = -11824
$ws.Range("H61").Value = 1293.762
$ws.Range("I61").Value = 1201.0968
$ws.Range("J61").Value = 1554.909
$ws.Range("K61").Value = 1201.0968
$ws.Range("L61").Value = 1554.909
$ws.Range("M61").Value = -989.0968
$ws.Range("N61").Value = -1978.909
$ws.Range("H97").Value = 498.1
$ws.Range("I97").Value = 425.66666
$ws.Range("K97").Value = 425.66666
$ws.Range("M97").Value = 70.33334000000002
$ws.Range("H102").Value = 18519874
$ws.Range("I102").Value = 18519874
$ws.Range("K102").Value = 18519874
$ws.Range("M102").Value = -18518252
$ws.Range("H110").Value = 1635.15
$ws.Range("I110").Value = 1259.3334
$ws.Range("K110").Value = 1259.3334
$ws.Range("M110").Value = 785.6666
$ws.Range("H132").Value = 1409.4814
$ws.Range("I132").Value = 1137.973
$ws.Range("J132").Value = 2000.4117
$ws.Range("K132").Value = 3413.919
$ws.Range("L132").Value = 6001.2351
$ws.Range("M132").Value = -883.9189999999999
$ws.Range("N132").Value = -11061.2351
$ws.Range("H136").Value = 1293.762
$ws.Range("I136").Value = 1201.0968
$ws.Range("J136").Value = 1554.909
$ws.Range("K136").Value = 3603.2904
$ws.Range("L136").Value = 4664.727000000001
$ws.Range("M136").Value = -1053.2904
$ws.Range("N136").Value = -9764.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 9616214
$ws.Range("I94").Value = 13158744
$ws.Range("J94").Value = 775.7143
$ws.Range("K94").Value = 13158744
$ws.Range("L94").Value = 775.7143
$ws.Range("M94").Value = -13158293
$ws.Range("N94").Value = -1677.7143
$ws.Range("H134").Value = 3697.7173
$ws.Range("I134").Value = 907.0263
$ws.Range("J134").Value = 16953.5
$ws.Range("K134").Value = 2721.0789
$ws.Range("L134").Value = 50860.5
$ws.Range("M134").Value = -186.0789
$ws.Range("N134").Value = -55930.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 47055.934
$ws.Range("I22").Value = 416.125
$ws.Range("J22").Value = 100358.57
$ws.Range("K22").Value = 416.125
$ws.Range("L22").Value = 100358.57
$ws.Range("M22").Value = -66.125
$ws.Range("N22").Value = -101058.57
$ws.Range("H31").Value = 1559.525
$ws.Range("I31").Value = 1578.2424
$ws.Range("K31").Value = 1578.2424
$ws.Range("M31").Value = -1283.2424
$ws.Range("H34").Value = 1559.525
$ws.Range("I34").Value = 1578.2424
$ws.Range("K34").Value = 1578.2424
$ws.Range("M34").Value = -1376.2424
$ws.Range("H50").Value = 19626.666
$ws.Range("J50").Value = 19626.666
$ws.Range("L50").Value = 19626.666
$ws.Range("N50").Value = -20876.666
$ws.Range("H58").Value = 788.55554
$ws.Range("I58").Value = 713.7406999999999
$ws.Range("J58").Value = 1013
$ws.Range("K58").Value = 713.7406999999999
$ws.Range("L58").Value = 1013
$ws.Range("M58").Value = -510.7406999999999
$ws.Range("N58").Value = -1419
$ws.Range("H99").Value = 2287.7778
$ws.Range("I99").Value = 2287.7778
$ws.Range("K99").Value = 2287.7778
$ws.Range("M99").Value = -789.7777999999998
$ws.Range("H122").Value = 1315.138
$ws.Range("I122").Value = 956.5833
$ws.Range("J122").Value = 1568.2354
$ws.Range("K122").Value = 2869.7499
$ws.Range("L122").Value = 4704.706200000001
$ws.Range("M122").Value = -419.7498999999998
$ws.Range("N122").Value = -9604.706200000001
$ws.Range("H126").Value = 2287.7778
$ws.Range("I126").Value = 2287.7778
$ws.Range("K126").Value = 6863.3334
$ws.Range("M126").Value = -4393.3334
$ws.Range("H132").Value = 3782.7917
$ws.Range("I132").Value = 4157.7837
$ws.Range("J132").Value = 2521.4546
$ws.Range("K132").Value = 12473.3511
$ws.Range("L132").Value = 7564.3638
$ws.Range("M132").Value = -9943.3511
$ws.Range("N132").Value = -12624.3638
$ws.Range("H136").Value = 788.55554
$ws.Range("I136").Value = 713.7406999999999
$ws.Range("J136").Value = 1013
$ws.Range("K136").Value = 2141.2221
$ws.Range("L136").Value = 3039
$ws.Range("M136").Value = 408.7779
$ws.Range("N136").Value = -8139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 672.38464
$ws.Range("I23").Value = 850
$ws.Range("J23").Value = 561.375
$ws.Range("K23").Value = 2550
$ws.Range("L23").Value = 1684.125
$ws.Range("M23").Value = -2315
$ws.Range("N23").Value = -2154.125
$ws.Range("H61").Value = 153.33333
$ws.Range("I61").Value = 107.77778
$ws.Range("J61").Value = 290
$ws.Range("K61").Value = 323.33334
$ws.Range("L61").Value = 870
$ws.Range("M61").Value = -108.33334
$ws.Range("N61").Value = -1300
$ws.Range("H68").Value = 1386.0714
$ws.Range("I68").Value = 1150.3334
$ws.Range("J68").Value = 1562.875
$ws.Range("K68").Value = 3451.0002
$ws.Range("L68").Value = 4688.625
$ws.Range("M68").Value = -2640.0002
$ws.Range("N68").Value = -6310.625
$ws.Range("H71").Value = 1386.0714
$ws.Range("I71").Value = 1150.3334
$ws.Range("J71").Value = 1562.875
$ws.Range("K71").Value = 10353.0006
$ws.Range("L71").Value = 14065.875
$ws.Range("M71").Value = -6297.000599999999
$ws.Range("N71").Value = -22177.875
$ws.Range("H76").Value = 6984.4165
$ws.Range("J76").Value = 6709.091
$ws.Range("L76").Value = 20127.273
$ws.Range("N76").Value = -20893.273
$ws.Range("H79").Value = 6984.4165
$ws.Range("J79").Value = 6709.091
$ws.Range("L79").Value = 20127.273
$ws.Range("N79").Value = -22779.273
$ws.Range("H132").Value = 1494.7273
$ws.Range("I132").Value = 1206.2858
$ws.Range("K132").Value = 10856.5722
$ws.Range("M132").Value = -8326.572200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 622.55
$ws.Range("I107").Value = 779.1818
$ws.Range("J107").Value = 431.1111
$ws.Range("K107").Value = 779.1818
$ws.Range("L107").Value = 431.1111
$ws.Range("M107").Value = 1140.8182
$ws.Range("N107").Value = -4271.1111
$ws.Range("H113").Value = 1402.3529
$ws.Range("I113").Value = 1362.6428
$ws.Range("J113").Value = 1587.6666
$ws.Range("K113").Value = 1362.6428
$ws.Range("L113").Value = 1587.6666
$ws.Range("M113").Value = 807.3571999999999
$ws.Range("N113").Value = -5927.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 611.0909
$ws.Range("I93").Value = 577.375
$ws.Range("J93").Value = 701
$ws.Range("K93").Value = 577.375
$ws.Range("L93").Value = 701
$ws.Range("M93").Value = 670.625
$ws.Range("N93").Value = -3197
$ws.Range("H100").Value = 1155.3334
$ws.Range("I100").Value = 729
$ws.Range("J100").Value = 1688.25
$ws.Range("K100").Value = 729
$ws.Range("L100").Value = 1688.25
$ws.Range("M100").Value = -188
$ws.Range("N100").Value = -2770.25
$ws.Range("H122").Value = 41686150
$ws.Range("I122").Value = 41686150
$ws.Range("K122").Value = 125058450
$ws.Range("M122").Value = -125056000
$ws.Range("H132").Value = 21960.18
$ws.Range("I132").Value = 1496.1923
$ws.Range("J132").Value = 44129.5
$ws.Range("K132").Value = 4488.5769
$ws.Range("L132").Value = 132388.5
$ws.Range("M132").Value = -1958.5769
$ws.Range("N132").Value = -137448.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3415.4546
$ws.Range("I96").Value = 3575
$ws.Range("J96").Value = 2990
$ws.Range("K96").Value = 3575
$ws.Range("L96").Value = 2990
$ws.Range("M96").Value = -2202
$ws.Range("N96").Value = -5736
$ws.Range("H104").Value = 25030
$ws.Range("J104").Value = 25030
$ws.Range("L104").Value = 25030
$ws.Range("N104").Value = -32018
$ws.Range("H132").Value = 1971.9615
$ws.Range("I132").Value = 2076.4524
$ws.Range("J132").Value = 1533.1
$ws.Range("K132").Value = 6229.3572
$ws.Range("L132").Value = 4599.299999999999
$ws.Range("M132").Value = -3699.3572
$ws.Range("N132").Value = -9659.299999999999
$ws.Range("H133").Value = 46000
$ws.Range("J133").Value = 46000
$ws.Range("L133").Value = 46000
$ws.Range("N133").Value = -56120
$ws.Range("H135").Value = 49886
$ws.Range("J135").Value = 49886
$ws.Range("L135").Value = 49886
$ws.Range("N135").Value = -60026
$ws.Range("H136").Value = 506.73334
$ws.Range("I136").Value = 338.04544
$ws.Range("J136").Value = 970.625
$ws.Range("K136").Value = 1014.13632
$ws.Range("L136").Value = 2911.875
$ws.Range("M136").Value = 1535.86368
$ws.Range("N136").Value = -8011.875
